# This workbook ships with a protected worksheet, so cell values cannot be
# written until the sheet is unprotected. We unprotect, apply the data
# updates from the new holdings snapshot, then re-protect the sheet so it
# ends up in the same (protected) state it started in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer note (A38).
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-04 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) figures for each holding row.
$ws.Range("D2").Value = 0.03617999716608691
$ws.Range("E2").Value = -0.003869969040247723
$ws.Range("D3").Value = 0.02038889983838142
$ws.Range("E3").Value = 0.002337358784573418
$ws.Range("D4").Value = 0.01915798566516538
$ws.Range("E4").Value = -0.0006059381943042119
$ws.Range("D5").Value = 0.03780723148838938
$ws.Range("E5").Value = -0.0003501400560222967
$ws.Range("D6").Value = 0.0342675386049378
$ws.Range("E6").Value = -0.0007993605115907965
$ws.Range("D7").Value = 0.01977654854944738
$ws.Range("E7").Value = -0.001351612280362957
$ws.Range("D8").Value = 0.03709137058189461
$ws.Range("E8").Value = -0.006170212765957417
$ws.Range("D9").Value = 0.02027953502848449
$ws.Range("E9").Value = 0.0009038322487346928
$ws.Range("D10").Value = 0.025638105225138
$ws.Range("E10").Value = 0.0006950650382284707
$ws.Range("D11").Value = 0.02368807160360531
$ws.Range("E11").Value = -0.003761418592154731
$ws.Range("D12").Value = 0.05700060288115176
$ws.Range("E12").Value = -0.002861912711662251
$ws.Range("D13").Value = 0.02502865607501408
$ws.Range("E13").Value = 0.0003661662394727205
$ws.Range("D14").Value = 0.02757591931483863
$ws.Range("E14").Value = 0.005601816805450577
$ws.Range("D15").Value = 0.03302857990659899
$ws.Range("E15").Value = 0.01174033149171261
$ws.Range("D16").Value = 0.019670951428532
$ws.Range("E16").Value = 0.003012804418779824
$ws.Range("D17").Value = 0.03054628384434822
$ws.Range("E17").Value = -0.005667138928243998
$ws.Range("D18").Value = 0.04210036033871932
$ws.Range("E18").Value = -0.00022977941176483
$ws.Range("D19").Value = 0.1261553630990606
$ws.Range("E19").Value = 0
$ws.Range("D20").Value = 0.008897550274081555
$ws.Range("E20").Value = 0.0001316135825215614
$ws.Range("D21").Value = 0.01559334457089001
$ws.Range("E21").Value = -0.00329127811300034
$ws.Range("D22").Value = 0.01658826903559942
$ws.Range("E22").Value = -0.006574484737803421
$ws.Range("D23").Value = 0.01636607721511982
$ws.Range("E23").Value = -0.01186842997626303
$ws.Range("D24").Value = 0.02168096258547372
$ws.Range("E24").Value = -0.009088155104513773
$ws.Range("D25").Value = 0.01220216991391391
$ws.Range("E25").Value = 0.00475258596589323
$ws.Range("D26").Value = 0.04197867416755839
$ws.Range("E26").Value = -0.006167448968453226
$ws.Range("D27").Value = 0.02387279019313327
$ws.Range("E27").Value = 0
$ws.Range("D28").Value = 0.04565675324900804
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0.05548482087210086
$ws.Range("E29").Value = -0.007374100719424459
$ws.Range("D30").Value = 0.01333700801810038
$ws.Range("E30").Value = -0.01394169835234471
$ws.Range("D31").Value = 0.0205738730014892
$ws.Range("E31").Value = 0.003452243958573131
$ws.Range("D32").Value = 0.01388062444047909
$ws.Range("E32").Value = -0.01100412654745531
$ws.Range("D33").Value = 0.04182470806647446
$ws.Range("E33").Value = -0.001029866117404854
$ws.Range("D34").Value = 0.01668037375678356
$ws.Range("E34").Value = -0.01675142087944959
$ws.Range("E35").Value = -0.001939086957599701

$ws.Protect()
